$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 4500
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 4500
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -3564
$ws.Range("N74").Value = -10872
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 4500
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 22500
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -17820
$ws.Range("N77").Value = -54360
$ws.Range("H98").Value = 35576896
$ws.Range("I98").Value = 15386509
$ws.Range("K98").Value = 15386509
$ws.Range("M98").Value = -15385011
$ws.Range("H100").Value = 15386498
$ws.Range("I100").Value = 1138.8572
$ws.Range("J100").Value = 33336084
$ws.Range("K100").Value = 1138.8572
$ws.Range("L100").Value = 33336084
$ws.Range("M100").Value = -597.8571999999999
$ws.Range("N100").Value = -33337166
$ws.Range("H113").Value = 4168429.5
$ws.Range("I113").Value = 5883959.5
$ws.Range("J113").Value = 2142.8572
$ws.Range("K113").Value = 5883959.5
$ws.Range("L113").Value = 2142.8572
$ws.Range("M113").Value = -5880705.5
$ws.Range("N113").Value = -8650.8572
$ws.Range("H116").Value = 13731896
$ws.Range("J116").Value = 18528542
$ws.Range("L116").Value = 18528542
$ws.Range("N116").Value = -18535426
$ws.Range("H122").Value = 35576896
$ws.Range("I122").Value = 15386509
$ws.Range("K122").Value = 46159527
$ws.Range("M122").Value = -46157077
$ws.Range("H132").Value = 6174479
$ws.Range("I132").Value = 1736.1177
$ws.Range("J132").Value = 111111110
$ws.Range("K132").Value = 5208.3531
$ws.Range("L132").Value = 333333330
$ws.Range("M132").Value = -2678.3531
$ws.Range("N132").Value = -333338390
$ws.Range("H137").Value = 17454012
$ws.Range("I137").Value = 3472974
$ws.Range("J137").Value = 89356500
$ws.Range("K137").Value = 10418922
$ws.Range("L137").Value = 268069500
$ws.Range("M137").Value = -10416372
$ws.Range("N137").Value = -268074600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 633.8
$ws.Range("J97").Value = 417.25
$ws.Range("L97").Value = 417.25
$ws.Range("N97").Value = -1409.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1913.3077
$ws.Range("I94").Value = 1572.7142
$ws.Range("K94").Value = 1572.7142
$ws.Range("M94").Value = -1121.7142
$ws.Range("H105").Value = 2000
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2370567.5
$ws.Range("I31").Value = 1303558.1
$ws.Range("K31").Value = 1303558.1
$ws.Range("M31").Value = -1303263.1
$ws.Range("H34").Value = 2370567.5
$ws.Range("I34").Value = 1303558.1
$ws.Range("K34").Value = 1303558.1
$ws.Range("M34").Value = -1303356.1
$ws.Range("H62").Value = 2669.2307
$ws.Range("I62").Value = 2609.0908
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2609.0908
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1985.0908
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2669.2307
$ws.Range("I65").Value = 2609.0908
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 13045.454
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -9925.454
$ws.Range("N65").Value = -21240
$ws.Range("H105").Value = 8920.625
$ws.Range("I105").Value = 2175.7144
$ws.Range("J105").Value = 14166.667
$ws.Range("K105").Value = 2175.7144
$ws.Range("L105").Value = 14166.667
$ws.Range("M105").Value = -428.7143999999998
$ws.Range("N105").Value = -17660.667
$ws.Range("H107").Value = 586.95
$ws.Range("I107").Value = 257
$ws.Range("J107").Value = 856.9091
$ws.Range("K107").Value = 257
$ws.Range("L107").Value = 856.9091
$ws.Range("M107").Value = 1663
$ws.Range("N107").Value = -4696.9091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2878139.2
$ws.Range("I70").Value = 1742033.4
$ws.Range("J70").Value = 4791580.5
$ws.Range("K70").Value = 1742033.4
$ws.Range("L70").Value = 4791580.5
$ws.Range("M70").Value = -1741763.4
$ws.Range("N70").Value = -4792120.5
$ws.Range("H73").Value = 2878139.2
$ws.Range("I73").Value = 1742033.4
$ws.Range("J73").Value = 4791580.5
$ws.Range("K73").Value = 1742033.4
$ws.Range("L73").Value = 4791580.5
$ws.Range("M73").Value = -1741097.4
$ws.Range("N73").Value = -4793452.5
$ws.Range("H102").Value = 4820.769
$ws.Range("I102").Value = 5175.6523
$ws.Range("K102").Value = 5175.6523
$ws.Range("M102").Value = -3553.6523
$ws.Range("H107").Value = 316.5
$ws.Range("I107").Value = 88.666664
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 88.666664
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1831.333336
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1442.0476
$ws.Range("I7").Value = 1050.3
$ws.Range("K7").Value = 1050.3
$ws.Range("M7").Value = -938.3
$ws.Range("H100").Value = 1880.2632
$ws.Range("I100").Value = 1558
$ws.Range("J100").Value = 1995.3572
$ws.Range("K100").Value = 1558
$ws.Range("L100").Value = 1995.3572
$ws.Range("M100").Value = -1017
$ws.Range("N100").Value = -3077.3572
$ws.Range("H126").Value = 1442.0476
$ws.Range("I126").Value = 1050.3
$ws.Range("K126").Value = 3150.9
$ws.Range("M126").Value = -680.8999999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1090.1282
$ws.Range("I122").Value = 971.44446
$ws.Range("J122").Value = 1357.1666
$ws.Range("K122").Value = 2914.33338
$ws.Range("L122").Value = 4071.4998
$ws.Range("M122").Value = -464.33338
$ws.Range("N122").Value = -8971.4998
$ws.Range("H126").Value = 19232590
$ws.Range("I126").Value = 22727932
$ws.Range("K126").Value = 68183796
$ws.Range("M126").Value = -68181326
$ws.Range("H132").Value = 1352110
$ws.Range("I132").Value = 987585.8
$ws.Range("J132").Value = 2333521.5
$ws.Range("K132").Value = 2962757.4
$ws.Range("L132").Value = 7000564.5
$ws.Range("M132").Value = -2960227.4
$ws.Range("N132").Value = -7005624.5
$ws.Range("H136").Value = 10474.708
$ws.Range("I136").Value = 7980.25
$ws.Range("J136").Value = 15463.625
$ws.Range("K136").Value = 23940.75
$ws.Range("L136").Value = 46390.875
$ws.Range("M136").Value = -21390.75
$ws.Range("N136").Value = -51490.875
